# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F value for sheet "展览" (rId1 / sheet1.xml)
$updatesDisplay = @{
    3  = 2933
    4  = 190
    5  = 152
    7  = 1585
    9  = 77
    10 = 25
    11 = 1308
    13 = 430
    17 = 117
    19 = 96
    20 = 2989
    21 = 360
    22 = 23
    23 = 17
    24 = 78
}

# Row -> new F value for sheet "全部类型" (rId4 / sheet4.xml)
# Identical to the above except row 11 (1309 instead of 1308).
$updatesAll = @{
    3  = 2933
    4  = 190
    5  = 152
    7  = 1585
    9  = 77
    10 = 25
    11 = 1309
    13 = 430
    17 = 117
    19 = 96
    20 = 2989
    21 = 360
    22 = 23
    23 = 17
    24 = 78
}

$wsDisplay = $wb.Worksheets.Item("展览")
foreach ($row in $updatesDisplay.Keys) {
    $wsDisplay.Cells.Item($row, 6).Value = $updatesDisplay[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $updatesAll[$row]
}
